$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Montenegro Prva Liga")

# Row 17
$ws.Range("B17").Value = 6815303
$ws.Range("F17").Value = "FK Mornar Bar"
$ws.Range("G17").Value = "FK Arsenal"
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = "A"
$ws.Range("K17").Value = 2.4
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 2.75
$ws.Range("N17").Value = 2.4
$ws.Range("O17").Value = 3
$ws.Range("P17").Value = 2.75
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 1.775
$ws.Range("S17").Value = 2.025
$ws.Range("T17").Value = 2
$ws.Range("U17").Value = 1.85
$ws.Range("V17").Value = 1.95
$ws.Range("W17").Value = -1
$ws.Range("X17").Value = -1
$ws.Range("Y17").Value = 1.75
$ws.Range("Z17").Value = -1
$ws.Range("AA17").Value = 1.025
$ws.Range("AB17").Value = -1
$ws.Range("AC17").Value = 0.95
# Row 18
$ws.Range("B18").Value = 6815306
$ws.Range("F18").Value = "OFK Petrovac"
$ws.Range("G18").Value = "FK Rudar Pljevlja"
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = "H"
$ws.Range("K18").Value = 2.1
$ws.Range("L18").Value = 3.1
$ws.Range("M18").Value = 3.2
$ws.Range("N18").Value = 1.615
$ws.Range("O18").Value = 3.5
$ws.Range("P18").Value = 5
$ws.Range("Q18").Value = -0.75
$ws.Range("R18").Value = 1.85
$ws.Range("S18").Value = 1.95
$ws.Range("T18").Value = 2.25
$ws.Range("U18").Value = 1.85
$ws.Range("V18").Value = 1.95
$ws.Range("W18").Value = 0.615
$ws.Range("X18").Value = -1
$ws.Range("Y18").Value = -1
$ws.Range("Z18").Value = 0.8500000000000001
$ws.Range("AA18").Value = -1
$ws.Range("AB18").Value = -0.5
$ws.Range("AC18").Value = 0.475
# Row 19
$ws.Range("B19").Value = 6815304
$ws.Range("F19").Value = "FK Jedinstvo Bijelo Polje"
$ws.Range("G19").Value = "Sutjeska Niksic"
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = "D"
$ws.Range("K19").Value = 5.5
$ws.Range("L19").Value = 3.75
$ws.Range("M19").Value = 1.5
$ws.Range("N19").Value = 3.6
$ws.Range("O19").Value = 3.2
$ws.Range("P19").Value = 1.909
$ws.Range("Q19").Value = 0.5
$ws.Range("R19").Value = 1.825
$ws.Range("S19").Value = 1.975
$ws.Range("T19").Value = 2.25
$ws.Range("U19").Value = 1.875
$ws.Range("V19").Value = 1.925
$ws.Range("W19").Value = -1
$ws.Range("X19").Value = 2.2
$ws.Range("Y19").Value = -1
$ws.Range("Z19").Value = 0.825
$ws.Range("AA19").Value = -1
$ws.Range("AB19").Value = -1
$ws.Range("AC19").Value = 0.925
# Row 20
$ws.Range("B20").Value = 6815305
$ws.Range("F20").Value = "Buducnost Podgorica"
$ws.Range("G20").Value = "FK Jezero"
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = "D"
$ws.Range("K20").Value = 1.3
$ws.Range("L20").Value = 5
$ws.Range("M20").Value = 7
$ws.Range("N20").Value = 1.571
$ws.Range("O20").Value = 4
$ws.Range("P20").Value = 4.2
$ws.Range("Q20").Value = -0.75
$ws.Range("R20").Value = 1.75
$ws.Range("S20").Value = 1.95
$ws.Range("T20").Value = 2.5
$ws.Range("U20").Value = 1.95
$ws.Range("V20").Value = 1.85
$ws.Range("W20").Value = -1
$ws.Range("X20").Value = 3
$ws.Range("Y20").Value = -1
$ws.Range("Z20").Value = -1
$ws.Range("AA20").Value = 0.95
$ws.Range("AB20").Value = -1
$ws.Range("AC20").Value = 0.8500000000000001
# Row 21
$ws.Range("B21").Value = 6815422
$ws.Range("F21").Value = "OFK Mladost DG"
$ws.Range("G21").Value = "FK Decic Tuzi"
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = "A"
$ws.Range("K21").Value = 2.4
$ws.Range("L21").Value = 3
$ws.Range("M21").Value = 2.75
$ws.Range("N21").Value = 3.1
$ws.Range("O21").Value = 3
$ws.Range("P21").Value = 2.15
$ws.Range("Q21").Value = 0.25
$ws.Range("R21").Value = 1.875
$ws.Range("S21").Value = 1.925
$ws.Range("T21").Value = 2.25
$ws.Range("U21").Value = 2.025
$ws.Range("V21").Value = 1.775
$ws.Range("W21").Value = -1
$ws.Range("X21").Value = -1
$ws.Range("Y21").Value = 1.15
$ws.Range("Z21").Value = -1
$ws.Range("AA21").Value = 0.925
$ws.Range("AB21").Value = -0.5
$ws.Range("AC21").Value = 0.3875
# Row 38
$ws.Range("B38").Value = 6815322
$ws.Range("F38").Value = "OFK Mladost DG"
$ws.Range("G38").Value = "FK Arsenal"
$ws.Range("H38").Value = 2
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = "H"
$ws.Range("K38").Value = 2.375
$ws.Range("L38").Value = 3
$ws.Range("M38").Value = 2.75
$ws.Range("N38").Value = 2.625
$ws.Range("O38").Value = 3
$ws.Range("P38").Value = 2.4
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = 2
$ws.Range("S38").Value = 1.8
$ws.Range("T38").Value = 2
$ws.Range("U38").Value = 1.725
$ws.Range("V38").Value = 1.975
$ws.Range("W38").Value = 1.625
$ws.Range("X38").Value = -1
$ws.Range("Y38").Value = -1
$ws.Range("Z38").Value = 1
$ws.Range("AA38").Value = -1
$ws.Range("AB38").Value = 0.7250000000000001
$ws.Range("AC38").Value = -1
# Row 39
$ws.Range("B39").Value = 6815321
$ws.Range("F39").Value = "OFK Petrovac"
$ws.Range("G39").Value = "FK Jedinstvo Bijelo Polje"
$ws.Range("H39").Value = 1
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = "D"
$ws.Range("K39").Value = 1.8
$ws.Range("L39").Value = 3.4
$ws.Range("M39").Value = 3.75
$ws.Range("N39").Value = 1.6
$ws.Range("O39").Value = 3.5
$ws.Range("P39").Value = 4.75
$ws.Range("Q39").Value = -0.75
$ws.Range("R39").Value = 1.825
$ws.Range("S39").Value = 1.975
$ws.Range("T39").Value = 2.5
$ws.Range("U39").Value = 1.95
$ws.Range("V39").Value = 1.75
$ws.Range("W39").Value = -1
$ws.Range("X39").Value = 2.5
$ws.Range("Y39").Value = -1
$ws.Range("Z39").Value = -1
$ws.Range("AA39").Value = 0.9750000000000001
$ws.Range("AB39").Value = -1
$ws.Range("AC39").Value = 0.75
# Row 42
$ws.Range("B42").Value = 7246029
$ws.Range("F42").Value = "Sutjeska Niksic"
$ws.Range("G42").Value = "FK Rudar Pljevlja"
$ws.Range("H42").Value = 1
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = "H"
$ws.Range("K42").Value = 1.3
$ws.Range("L42").Value = 4.333
$ws.Range("M42").Value = 9
$ws.Range("N42").Value = 1.363
$ws.Range("O42").Value = 4.333
$ws.Range("P42").Value = 7.5
$ws.Range("Q42").Value = -1.5
$ws.Range("R42").Value = 1.95
$ws.Range("S42").Value = 1.85
$ws.Range("T42").Value = 3
$ws.Range("U42").Value = 1.975
$ws.Range("V42").Value = 1.825
$ws.Range("W42").Value = 0.363
$ws.Range("X42").Value = -1
$ws.Range("Y42").Value = -1
$ws.Range("Z42").Value = -1
$ws.Range("AA42").Value = 0.8500000000000001
$ws.Range("AB42").Value = -1
$ws.Range("AC42").Value = 0.825
# Row 43
$ws.Range("B43").Value = 6815327
$ws.Range("F43").Value = "FK Arsenal"
$ws.Range("G43").Value = "FK Decic Tuzi"
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = "A"
$ws.Range("K43").Value = 3
$ws.Range("L43").Value = 3
$ws.Range("M43").Value = 2.25
$ws.Range("N43").Value = 3.25
$ws.Range("O43").Value = 3
$ws.Range("P43").Value = 2.1
$ws.Range("Q43").Value = 0.25
$ws.Range("R43").Value = 1.925
$ws.Range("S43").Value = 1.875
$ws.Range("T43").Value = 2
$ws.Range("U43").Value = 1.775
$ws.Range("V43").Value = 2.025
$ws.Range("W43").Value = -1
$ws.Range("X43").Value = -1
$ws.Range("Y43").Value = 1.1
$ws.Range("Z43").Value = -1
$ws.Range("AA43").Value = 0.875
$ws.Range("AB43").Value = -1
$ws.Range("AC43").Value = 1.025
# Row 53
$ws.Range("B53").Value = 6815334
$ws.Range("F53").Value = "Sutjeska Niksic"
$ws.Range("G53").Value = "FK Mornar Bar"
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = "A"
$ws.Range("K53").Value = 1.444
$ws.Range("L53").Value = 4
$ws.Range("M53").Value = 6.5
$ws.Range("N53").Value = 1.444
$ws.Range("O53").Value = 4
$ws.Range("P53").Value = 6.5
$ws.Range("Q53").Value = -1.25
$ws.Range("R53").Value = 2
$ws.Range("S53").Value = 1.8
$ws.Range("T53").Value = 2.5
$ws.Range("U53").Value = 2
$ws.Range("V53").Value = 1.8
$ws.Range("W53").Value = -1
$ws.Range("X53").Value = -1
$ws.Range("Y53").Value = 5.5
$ws.Range("Z53").Value = -1
$ws.Range("AA53").Value = 0.8
$ws.Range("AB53").Value = -1
$ws.Range("AC53").Value = 0.8
# Row 54
$ws.Range("B54").Value = 7279987
$ws.Range("F54").Value = "FK Jezero"
$ws.Range("G54").Value = "FK Arsenal"
$ws.Range("H54").Value = 1
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = "D"
$ws.Range("K54").Value = 2.1
$ws.Range("L54").Value = 3
$ws.Range("M54").Value = 3.25
$ws.Range("N54").Value = 2.05
$ws.Range("O54").Value = 3
$ws.Range("P54").Value = 3.4
$ws.Range("Q54").Value = -0.25
$ws.Range("R54").Value = 1.8
$ws.Range("S54").Value = 2
$ws.Range("T54").Value = 2
$ws.Range("U54").Value = 1.925
$ws.Range("V54").Value = 1.875
$ws.Range("W54").Value = -1
$ws.Range("X54").Value = 2
$ws.Range("Y54").Value = -1
$ws.Range("Z54").Value = -0.5
$ws.Range("AA54").Value = 0.5
$ws.Range("AB54").Value = 0
$ws.Range("AC54").Value = -0
# Row 59
$ws.Range("B59").Value = 6815427
$ws.Range("F59").Value = "FK Mornar Bar"
$ws.Range("G59").Value = "OFK Mladost DG"
$ws.Range("H59").Value = 2
$ws.Range("I59").Value = 1
$ws.Range("J59").Value = "H"
$ws.Range("K59").Value = 1.833
$ws.Range("L59").Value = 3.1
$ws.Range("M59").Value = 4
$ws.Range("N59").Value = 1.833
$ws.Range("O59").Value = 3.1
$ws.Range("P59").Value = 4
$ws.Range("Q59").Value = -0.5
$ws.Range("R59").Value = 1.875
$ws.Range("S59").Value = 1.925
$ws.Range("T59").Value = 2
$ws.Range("U59").Value = 1.775
$ws.Range("V59").Value = 2.025
$ws.Range("W59").Value = 0.833
$ws.Range("X59").Value = -1
$ws.Range("Y59").Value = -1
$ws.Range("Z59").Value = 0.875
$ws.Range("AA59").Value = -1
$ws.Range("AB59").Value = 0.7749999999999999
$ws.Range("AC59").Value = -1
# Row 60
$ws.Range("B60").Value = 6815338
$ws.Range("F60").Value = "OFK Petrovac"
$ws.Range("G60").Value = "FK Decic Tuzi"
$ws.Range("H60").Value = 2
$ws.Range("I60").Value = 3
$ws.Range("J60").Value = "A"
$ws.Range("K60").Value = 2.625
$ws.Range("L60").Value = 2.875
$ws.Range("M60").Value = 2.6
$ws.Range("N60").Value = 3.1
$ws.Range("O60").Value = 2.9
$ws.Range("P60").Value = 2.25
$ws.Range("Q60").Value = 0.25
$ws.Range("R60").Value = 1.8
$ws.Range("S60").Value = 2
$ws.Range("T60").Value = 2.25
$ws.Range("U60").Value = 1.975
$ws.Range("V60").Value = 1.725
$ws.Range("W60").Value = -1
$ws.Range("X60").Value = -1
$ws.Range("Y60").Value = 1.25
$ws.Range("Z60").Value = -1
$ws.Range("AA60").Value = 1
$ws.Range("AB60").Value = 0.9750000000000001
$ws.Range("AC60").Value = -1
# Row 62
$ws.Range("B62").Value = 7366683
$ws.Range("F62").Value = "FK Arsenal"
$ws.Range("G62").Value = "FK Mornar Bar"
$ws.Range("H62").Value = 2
$ws.Range("I62").Value = 2
$ws.Range("J62").Value = "D"
$ws.Range("K62").Value = 2.375
$ws.Range("L62").Value = 2.8
$ws.Range("M62").Value = 3
$ws.Range("N62").Value = 2.3
$ws.Range("O62").Value = 2.7
$ws.Range("P62").Value = 3.3
$ws.Range("Q62").Value = -0.25
$ws.Range("R62").Value = 2
$ws.Range("S62").Value = 1.8
$ws.Range("T62").Value = 1.75
$ws.Range("U62").Value = 1.875
$ws.Range("V62").Value = 1.925
$ws.Range("W62").Value = -1
$ws.Range("X62").Value = 1.7
$ws.Range("Y62").Value = -1
$ws.Range("Z62").Value = -0.5
$ws.Range("AA62").Value = 0.4
$ws.Range("AB62").Value = 0.875
$ws.Range("AC62").Value = -1
# Row 63
$ws.Range("B63").Value = 6815343
$ws.Range("F63").Value = "Sutjeska Niksic"
$ws.Range("G63").Value = "FK Jedinstvo Bijelo Polje"
$ws.Range("H63").Value = 2
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = "H"
$ws.Range("K63").Value = 1.333
$ws.Range("L63").Value = 4.2
$ws.Range("M63").Value = 8
$ws.Range("N63").Value = 1.333
$ws.Range("O63").Value = 4.2
$ws.Range("P63").Value = 8
$ws.Range("Q63").Value = -1.5
$ws.Range("R63").Value = 1.975
$ws.Range("S63").Value = 1.825
$ws.Range("T63").Value = 2.75
$ws.Range("U63").Value = 1.9
$ws.Range("V63").Value = 1.9
$ws.Range("W63").Value = 0.333
$ws.Range("X63").Value = -1
$ws.Range("Y63").Value = -1
$ws.Range("Z63").Value = 0.9750000000000001
$ws.Range("AA63").Value = -1
$ws.Range("AB63").Value = -1
$ws.Range("AC63").Value = 0.8999999999999999
# Row 75
$ws.Range("B75").Value = 6815358
$ws.Range("F75").Value = "OFK Petrovac"
$ws.Range("G75").Value = "FK Arsenal"
$ws.Range("H75").Value = 1
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = "D"
$ws.Range("K75").Value = 2.1
$ws.Range("L75").Value = 3.1
$ws.Range("M75").Value = 3.2
$ws.Range("N75").Value = 1.75
$ws.Range("O75").Value = 3.3
$ws.Range("P75").Value = 4.2
$ws.Range("Q75").Value = -0.5
$ws.Range("R75").Value = 1.8
$ws.Range("S75").Value = 2
$ws.Range("T75").Value = 2.25
$ws.Range("U75").Value = 1.95
$ws.Range("V75").Value = 1.85
$ws.Range("W75").Value = -1
$ws.Range("X75").Value = 2.3
$ws.Range("Y75").Value = -1
$ws.Range("Z75").Value = -1
$ws.Range("AA75").Value = 1
$ws.Range("AB75").Value = -0.5
$ws.Range("AC75").Value = 0.425
# Row 76
$ws.Range("B76").Value = 6815359
$ws.Range("F76").Value = "Buducnost Podgorica"
$ws.Range("G76").Value = "FK Jedinstvo Bijelo Polje"
$ws.Range("H76").Value = 3
$ws.Range("I76").Value = 2
$ws.Range("J76").Value = "H"
$ws.Range("K76").Value = 1.333
$ws.Range("L76").Value = 4.333
$ws.Range("M76").Value = 7.5
$ws.Range("N76").Value = 1.333
$ws.Range("O76").Value = 4.333
$ws.Range("P76").Value = 8
$ws.Range("Q76").Value = -1.5
$ws.Range("R76").Value = 1.875
$ws.Range("S76").Value = 1.925
$ws.Range("T76").Value = 2.75
$ws.Range("U76").Value = 1.8
$ws.Range("V76").Value = 2
$ws.Range("W76").Value = 0.333
$ws.Range("X76").Value = -1
$ws.Range("Y76").Value = -1
$ws.Range("Z76").Value = -1
$ws.Range("AA76").Value = 0.925
$ws.Range("AB76").Value = 0.8
$ws.Range("AC76").Value = -1
# Row 105
$ws.Range("B105").Value = 6815382
$ws.Range("F105").Value = "FK Jedinstvo Bijelo Polje"
$ws.Range("G105").Value = "Sutjeska Niksic"
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = "A"
$ws.Range("K105").Value = 3.75
$ws.Range("L105").Value = 3.3
$ws.Range("M105").Value = 1.833
$ws.Range("N105").Value = 4.75
$ws.Range("O105").Value = 3
$ws.Range("P105").Value = 1.75
$ws.Range("Q105").Value = 0.5
$ws.Range("R105").Value = 1.975
$ws.Range("S105").Value = 1.825
$ws.Range("T105").Value = 2
$ws.Range("U105").Value = 1.85
$ws.Range("V105").Value = 1.95
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = -1
$ws.Range("Y105").Value = 0.75
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = 0.825
$ws.Range("AB105").Value = -1
$ws.Range("AC105").Value = 0.95
# Row 106
$ws.Range("B106").Value = 6815434
$ws.Range("F106").Value = "OFK Mladost DG"
$ws.Range("G106").Value = "FK Decic Tuzi"
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 3
$ws.Range("J106").Value = "A"
$ws.Range("K106").Value = 4.6
$ws.Range("L106").Value = 3.6
$ws.Range("M106").Value = 1.615
$ws.Range("N106").Value = 8
$ws.Range("O106").Value = 4.75
$ws.Range("P106").Value = 1.25
$ws.Range("Q106").Value = 0.75
$ws.Range("R106").Value = 1.925
$ws.Range("S106").Value = 1.875
$ws.Range("T106").Value = 2.25
$ws.Range("U106").Value = 1.9
$ws.Range("V106").Value = 1.9
$ws.Range("W106").Value = -1
$ws.Range("X106").Value = -1
$ws.Range("Y106").Value = 0.25
$ws.Range("Z106").Value = -1
$ws.Range("AA106").Value = 0.875
$ws.Range("AB106").Value = 0.8999999999999999
$ws.Range("AC106").Value = -1
# Row 107
$ws.Range("B107").Value = 7890506
$ws.Range("F107").Value = "FK Mornar Bar"
$ws.Range("G107").Value = "FK Arsenal"
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = "D"
$ws.Range("K107").Value = 1.85
$ws.Range("L107").Value = 3.1
$ws.Range("M107").Value = 3.9
$ws.Range("N107").Value = 1.85
$ws.Range("O107").Value = 3.3
$ws.Range("P107").Value = 3.5
$ws.Range("Q107").Value = -0.5
$ws.Range("R107").Value = 1.925
$ws.Range("S107").Value = 1.875
$ws.Range("T107").Value = 2
$ws.Range("U107").Value = 1.95
$ws.Range("V107").Value = 1.85
$ws.Range("W107").Value = -1
$ws.Range("X107").Value = 2.3
$ws.Range("Y107").Value = -1
$ws.Range("Z107").Value = -1
$ws.Range("AA107").Value = 0.875
$ws.Range("AB107").Value = -1
$ws.Range("AC107").Value = 0.8500000000000001
# Row 108
$ws.Range("B108").Value = 7890508
$ws.Range("F108").Value = "OFK Petrovac"
$ws.Range("G108").Value = "FK Rudar Pljevlja"
$ws.Range("H108").Value = 1
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = "D"
$ws.Range("K108").Value = 1.75
$ws.Range("L108").Value = 3.1
$ws.Range("M108").Value = 4.5
$ws.Range("N108").Value = 1.8
$ws.Range("O108").Value = 3.2
$ws.Range("P108").Value = 4
$ws.Range("Q108").Value = -0.5
$ws.Range("R108").Value = 1.875
$ws.Range("S108").Value = 1.925
$ws.Range("T108").Value = 2.25
$ws.Range("U108").Value = 1.95
$ws.Range("V108").Value = 1.85
$ws.Range("W108").Value = -1
$ws.Range("X108").Value = 2.2
$ws.Range("Y108").Value = -1
$ws.Range("Z108").Value = -1
$ws.Range("AA108").Value = 0.925
$ws.Range("AB108").Value = -0.5
$ws.Range("AC108").Value = 0.425
# Row 110
$ws.Range("B110").Value = 7906320
$ws.Range("F110").Value = "FK Jezero"
$ws.Range("G110").Value = "OFK Petrovac"
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 2
$ws.Range("J110").Value = "A"
$ws.Range("K110").Value = 2.25
$ws.Range("L110").Value = 2.875
$ws.Range("M110").Value = 3.1
$ws.Range("N110").Value = 2.15
$ws.Range("O110").Value = 3
$ws.Range("P110").Value = 3.2
$ws.Range("Q110").Value = -0.25
$ws.Range("R110").Value = 1.9
$ws.Range("S110").Value = 1.9
$ws.Range("T110").Value = 1.75
$ws.Range("U110").Value = 1.85
$ws.Range("V110").Value = 1.95
$ws.Range("W110").Value = -1
$ws.Range("X110").Value = -1
$ws.Range("Y110").Value = 2.2
$ws.Range("Z110").Value = -1
$ws.Range("AA110").Value = 0.8999999999999999
$ws.Range("AB110").Value = 0.425
$ws.Range("AC110").Value = -0.5
# Row 111
$ws.Range("B111").Value = 7906319
$ws.Range("F111").Value = "FK Rudar Pljevlja"
$ws.Range("G111").Value = "OFK Mladost DG"
$ws.Range("H111").Value = 1
$ws.Range("I111").Value = 2
$ws.Range("J111").Value = "A"
$ws.Range("K111").Value = 2.3
$ws.Range("L111").Value = 2.875
$ws.Range("M111").Value = 3
$ws.Range("N111").Value = 1.909
$ws.Range("O111").Value = 3
$ws.Range("P111").Value = 4
$ws.Range("Q111").Value = -0.5
$ws.Range("R111").Value = 1.975
$ws.Range("S111").Value = 1.825
$ws.Range("T111").Value = 2
$ws.Range("U111").Value = 1.75
$ws.Range("V111").Value = 2.05
$ws.Range("W111").Value = -1
$ws.Range("X111").Value = -1
$ws.Range("Y111").Value = 3
$ws.Range("Z111").Value = -1
$ws.Range("AA111").Value = 0.825
$ws.Range("AB111").Value = 0.75
$ws.Range("AC111").Value = -1
# Row 117
$ws.Range("B117").Value = 6815389
$ws.Range("F117").Value = "FK Mornar Bar"
$ws.Range("G117").Value = "FK Jedinstvo Bijelo Polje"
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = "D"
$ws.Range("K117").Value = 1.909
$ws.Range("L117").Value = 3.1
$ws.Range("M117").Value = 3.75
$ws.Range("N117").Value = 1.75
$ws.Range("O117").Value = 3.2
$ws.Range("P117").Value = 4.333
$ws.Range("Q117").Value = -0.5
$ws.Range("R117").Value = 1.8
$ws.Range("S117").Value = 2
$ws.Range("T117").Value = 2
$ws.Range("U117").Value = 1.9
$ws.Range("V117").Value = 1.9
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 2.2
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 1
$ws.Range("AB117").Value = -1
$ws.Range("AC117").Value = 0.8999999999999999
# Row 118
$ws.Range("B118").Value = 6815393
$ws.Range("F118").Value = "FK Decic Tuzi"
$ws.Range("G118").Value = "FK Rudar Pljevlja"
$ws.Range("H118").Value = 3
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = "H"
$ws.Range("K118").Value = 1.285
$ws.Range("L118").Value = 4.75
$ws.Range("M118").Value = 8
$ws.Range("N118").Value = 1.333
$ws.Range("O118").Value = 4.5
$ws.Range("P118").Value = 7
$ws.Range("Q118").Value = -1.25
$ws.Range("R118").Value = 1.75
$ws.Range("S118").Value = 1.95
$ws.Range("T118").Value = 2.25
$ws.Range("U118").Value = 1.75
$ws.Range("V118").Value = 1.95
$ws.Range("W118").Value = 0.333
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 0.75
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = 0.75
$ws.Range("AC118").Value = -1
# Row 125
$ws.Range("H125").Value = 1
$ws.Range("I125").Value = 2
$ws.Range("J125").Value = "A"
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = 0.909
$ws.Range("Z125").Value = -1
$ws.Range("AA125").Value = 0.8
$ws.Range("AB125").Value = 0.825
$ws.Range("AC125").Value = -1
# Row 126
$ws.Range("B126").Value = 6815401
$ws.Range("F126").Value = "FK Decic Tuzi"
$ws.Range("G126").Value = "Sutjeska Niksic"
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = "D"
$ws.Range("K126").Value = 2.55
$ws.Range("L126").Value = 3
$ws.Range("M126").Value = 2.6
$ws.Range("N126").Value = 2.1
$ws.Range("O126").Value = 3.1
$ws.Range("P126").Value = 3.3
$ws.Range("Q126").Value = -0.25
$ws.Range("R126").Value = 1.825
$ws.Range("S126").Value = 1.975
$ws.Range("T126").Value = 2
$ws.Range("U126").Value = 1.925
$ws.Range("V126").Value = 1.875
$ws.Range("W126").Value = -1
$ws.Range("X126").Value = 2.1
$ws.Range("Y126").Value = -1
$ws.Range("Z126").Value = -0.5
$ws.Range("AA126").Value = 0.4875
$ws.Range("AB126").Value = -1
$ws.Range("AC126").Value = 0.875
# Row 127
$ws.Range("B127").Value = 6815402
$ws.Range("F127").Value = "FK Rudar Pljevlja"
$ws.Range("G127").Value = "FK Jezero"
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 1
$ws.Range("J127").Value = "A"
$ws.Range("K127").Value = 2.8
$ws.Range("L127").Value = 3
$ws.Range("M127").Value = 2.375
$ws.Range("N127").Value = 2.45
$ws.Range("O127").Value = 2.9
$ws.Range("P127").Value = 2.75
$ws.Range("Q127").Value = 0
$ws.Range("R127").Value = 1.775
$ws.Range("S127").Value = 2.025
$ws.Range("T127").Value = 1.75
$ws.Range("U127").Value = 1.825
$ws.Range("V127").Value = 1.975
$ws.Range("W127").Value = -1
$ws.Range("X127").Value = -1
$ws.Range("Y127").Value = 1.75
$ws.Range("Z127").Value = -1
$ws.Range("AA127").Value = 1.025
$ws.Range("AB127").Value = -1
$ws.Range("AC127").Value = 0.9750000000000001
# Row 128
$ws.Range("H128").Value = 1
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = "H"
$ws.Range("W128").Value = 1.15
$ws.Range("X128").Value = -1
$ws.Range("Y128").Value = -1
$ws.Range("Z128").Value = 0.95
$ws.Range("AA128").Value = -1
$ws.Range("AB128").Value = -1
$ws.Range("AC128").Value = 1.05
# Row 129
$ws.Range("H129").Value = 1
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = "H"
$ws.Range("W129").Value = 2.6
$ws.Range("X129").Value = -1
$ws.Range("Y129").Value = -1
$ws.Range("Z129").Value = 0.825
$ws.Range("AA129").Value = -1
$ws.Range("AB129").Value = -1
$ws.Range("AC129").Value = 1
